$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: split the title paragraph - insert an empty paragraph before
#         "Tour of Data Mining Algorithms" (same centered/bold pPr).
# ---------------------------------------------------------------------
$d.Paragraphs(1).Range.InsertParagraphBefore()

# ---------------------------------------------------------------------
# Step 2: append the new X-means / kD-tree sentences onto the end of the
#         Apriori paragraph.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "and their occurrence counts.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "and their occurrence counts. The x-means algorithm is finishing up implementation with a visual representation of the clusters. The kD tree was hard to implement but was better than just using a list of points that are in the data set. The algorithm is very computation intensive measuring a lot of distances takes quite a bit of time. ",
    2) | Out-Null

Write-Output "--- after steps 1-2 ---"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output ("$i`: [" + $d.Paragraphs($i).Range.Text + "]")
}

# ---------------------------------------------------------------------
# Step 3: the old CLI paragraph (#8) text is replaced with the new
#         X-means narrative, followed by two blank paragraphs and then
#         the (typo-fixed) CLI paragraph - all four share the same
#         "ind firstLine=720" paragraph formatting the CLI paragraph
#         already had.
# ---------------------------------------------------------------------
$cliPara = $d.Paragraphs(8)
$xmeans = "X-Means started with k-means and that worked well then we implemented the kd tree to speed up the process of finding clusters. While testing k-means we used a 2-Gaussian distribution and used a k of 2 with the k means. This worked well but was hard to know for sure without visualizing it, but weka helped figure that out.   There was a small issue when figuring out when to reevaluate the centroids a few implementations sent the program into an infinite loop but limiting the number of loops helped. "
$cliFixed = "The command line interface for our project has also been implemented fully for Apriori.  The Apache Commons CLI library was used to create the command line interface.  When calling the program from the command line, the algorithm and input file must be specified, along with algorithm-specific options such as the minimum support for an itemset to be considered frequent by the Apriori algorithm.  Optional options may be specified, such as specifying a specific name for the output file, or specifying the delimiter used to separate attributes in the input file.  A help message is also provided when the help option is given, or when invalid or incomplete options are specified."
$cliPara.Range.Text = $xmeans + "`r`r`r" + $cliFixed

Write-Output "--- after step 3 ---"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output ("$i`: [" + $d.Paragraphs($i).Range.Text + "]")
}

Write-Output "--- full listing with pPr info ---"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    Write-Output ("$i`: ind=" + $p.Range.ParagraphFormat.FirstLineIndent + " color=" + $p.Range.Font.Color + " [" + $p.Range.Text + "]")
}
